# feat: add 2022-Q1 data
#
# Inserts a new "2022-Q1" sheet (21 -> 22 holdings) between "2021-Q4" and
# "总计", and updates the "总计" (totals) sheet with a new leading row for
# the 2022-Q1 quarter, shifting the older quarters down by one row.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) Create the "2022-Q1" sheet by duplicating "2021-Q4" (sheet index 3) so it
#    inherits the same header / column-A styling (style index 2), then drop it
#    in right before "总计" and rename it.
# ---------------------------------------------------------------------------
$template = $wb.Worksheets.Item(3)
$template.Copy($null, $template)
$ws = $wb.Worksheets.Item(4)
$ws.Name = "2022-Q1"

# The template only carried 20 holdings (rows 2..21); this quarter has 22
# (rows 2..23), so extend the styled column-A border down two more rows.
$ws.Range("A21").Copy()
$ws.Range("A22:A23").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2) Write the 2022-Q1 holdings. Numeric-looking text (fund codes that can
#    carry leading zeros, and the pre-formatted percentage/ratio strings) is
#    entered with a leading apostrophe so it is stored as text, then
#    ClearFormats() strips the resulting quote-prefix style so the cell keeps
#    the plain/default formatting (matching the other quarter sheets).
# ---------------------------------------------------------------------------
$ws.Range("A2").Value = 0
$ws.Range("B2").Value = "'002666"
$ws.Range("B2").ClearFormats()
$ws.Range("C2").Value = '前海开源沪港深创新成长灵活配置混合A'
$ws.Range("D2").Value = "'11.96"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'81.64"
$ws.Range("E2").ClearFormats()
$ws.Range("F2").Value = "'8.00"
$ws.Range("F2").ClearFormats()
$ws.Range("G2").Value = "'0.9568"
$ws.Range("G2").ClearFormats()
$ws.Range("H2").Value = 2

$ws.Range("A3").Value = 1
$ws.Range("B3").Value = "'003293"
$ws.Range("B3").ClearFormats()
$ws.Range("C3").Value = '易方达科瑞灵活配置混合'
$ws.Range("D3").Value = "'34.67"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'78.17"
$ws.Range("E3").ClearFormats()
$ws.Range("F3").Value = "'2.53"
$ws.Range("F3").ClearFormats()
$ws.Range("G3").Value = "'0.8772"
$ws.Range("G3").ClearFormats()
$ws.Range("H3").Value = 7

$ws.Range("A4").Value = 2
$ws.Range("B4").Value = "'506005"
$ws.Range("B4").ClearFormats()
$ws.Range("C4").Value = '博时科创板三年定期开放混合'
$ws.Range("D4").Value = "'22.84"
$ws.Range("D4").ClearFormats()
$ws.Range("E4").Value = "'96.44"
$ws.Range("E4").ClearFormats()
$ws.Range("F4").Value = "'3.33"
$ws.Range("F4").ClearFormats()
$ws.Range("G4").Value = "'0.7606"
$ws.Range("G4").ClearFormats()
$ws.Range("H4").Value = 10

$ws.Range("A5").Value = 3
$ws.Range("B5").Value = "'010389"
$ws.Range("B5").ClearFormats()
$ws.Range("C5").Value = '易方达科益混合A'
$ws.Range("D5").Value = "'7.10"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'92.94"
$ws.Range("E5").ClearFormats()
$ws.Range("F5").Value = "'8.85"
$ws.Range("F5").ClearFormats()
$ws.Range("G5").Value = "'0.6284"
$ws.Range("G5").ClearFormats()
$ws.Range("H5").Value = 2

$ws.Range("A6").Value = 4
$ws.Range("B6").Value = "'011826"
$ws.Range("B6").ClearFormats()
$ws.Range("C6").Value = '汇添富健康生活一年持有期混合A'
$ws.Range("D6").Value = "'14.08"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'88.44"
$ws.Range("E6").ClearFormats()
$ws.Range("F6").Value = "'3.38"
$ws.Range("F6").ClearFormats()
$ws.Range("G6").Value = "'0.4759"
$ws.Range("G6").ClearFormats()
$ws.Range("H6").Value = 9

$ws.Range("A7").Value = 5
$ws.Range("B7").Value = "'011649"
$ws.Range("B7").ClearFormats()
$ws.Range("C7").Value = '易方达逆向投资混合A'
$ws.Range("D7").Value = "'7.49"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "'85.02"
$ws.Range("E7").ClearFormats()
$ws.Range("F7").Value = "'5.82"
$ws.Range("F7").ClearFormats()
$ws.Range("G7").Value = "'0.4359"
$ws.Range("G7").ClearFormats()
$ws.Range("H7").Value = 1

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = "'519019"
$ws.Range("B8").ClearFormats()
$ws.Range("C8").Value = '大成景阳领先混合'
$ws.Range("D8").Value = "'10.33"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'92.80"
$ws.Range("E8").ClearFormats()
$ws.Range("F8").Value = "'3.96"
$ws.Range("F8").ClearFormats()
$ws.Range("G8").Value = "'0.4091"
$ws.Range("G8").ClearFormats()
$ws.Range("H8").Value = 9

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = "'110012"
$ws.Range("B9").ClearFormats()
$ws.Range("C9").Value = '易方达科汇灵活配置混合'
$ws.Range("D9").Value = "'15.73"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'75.64"
$ws.Range("E9").ClearFormats()
$ws.Range("F9").Value = "'2.59"
$ws.Range("F9").ClearFormats()
$ws.Range("G9").Value = "'0.4074"
$ws.Range("G9").ClearFormats()
$ws.Range("H9").Value = 7

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = "'000124"
$ws.Range("B10").ClearFormats()
$ws.Range("C10").Value = '华宝服务优选混合'
$ws.Range("D10").Value = "'6.61"
$ws.Range("D10").ClearFormats()
$ws.Range("E10").Value = "'89.07"
$ws.Range("E10").ClearFormats()
$ws.Range("F10").Value = "'5.73"
$ws.Range("F10").ClearFormats()
$ws.Range("G10").Value = "'0.3788"
$ws.Range("G10").ClearFormats()
$ws.Range("H10").Value = 3

$ws.Range("A11").Value = 9
$ws.Range("B11").Value = "'240001"
$ws.Range("B11").ClearFormats()
$ws.Range("C11").Value = '华宝宝康消费品混合'
$ws.Range("D11").Value = "'11.24"
$ws.Range("D11").ClearFormats()
$ws.Range("E11").Value = "'62.57"
$ws.Range("E11").ClearFormats()
$ws.Range("F11").Value = "'2.97"
$ws.Range("F11").ClearFormats()
$ws.Range("G11").Value = "'0.3338"
$ws.Range("G11").ClearFormats()
$ws.Range("H11").Value = 8

$ws.Range("A12").Value = 10
$ws.Range("B12").Value = "'002667"
$ws.Range("B12").ClearFormats()
$ws.Range("C12").Value = '前海开源沪港深创新成长灵活配置混合C'
$ws.Range("D12").Value = "'3.25"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "'81.64"
$ws.Range("E12").ClearFormats()
$ws.Range("F12").Value = "'8.00"
$ws.Range("F12").ClearFormats()
$ws.Range("G12").Value = "'0.2600"
$ws.Range("G12").ClearFormats()
$ws.Range("H12").Value = 2

$ws.Range("A13").Value = 11
$ws.Range("B13").Value = "'001088"
$ws.Range("B13").ClearFormats()
$ws.Range("C13").Value = '华宝国策导向混合'
$ws.Range("D13").Value = "'3.45"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'85.75"
$ws.Range("E13").ClearFormats()
$ws.Range("F13").Value = "'6.75"
$ws.Range("F13").ClearFormats()
$ws.Range("G13").Value = "'0.2329"
$ws.Range("G13").ClearFormats()
$ws.Range("H13").Value = 2

$ws.Range("A14").Value = 12
$ws.Range("B14").Value = "'090016"
$ws.Range("B14").ClearFormats()
$ws.Range("C14").Value = '大成消费主题混合'
$ws.Range("D14").Value = "'4.23"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'93.78"
$ws.Range("E14").ClearFormats()
$ws.Range("F14").Value = "'4.28"
$ws.Range("F14").ClearFormats()
$ws.Range("G14").Value = "'0.1810"
$ws.Range("G14").ClearFormats()
$ws.Range("H14").Value = 8

$ws.Range("A15").Value = 13
$ws.Range("B15").Value = "'240002"
$ws.Range("B15").ClearFormats()
$ws.Range("C15").Value = '华宝宝康配置混合'
$ws.Range("D15").Value = "'4.54"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "'62.88"
$ws.Range("E15").ClearFormats()
$ws.Range("F15").Value = "'2.83"
$ws.Range("F15").ClearFormats()
$ws.Range("G15").Value = "'0.1285"
$ws.Range("G15").ClearFormats()
$ws.Range("H15").Value = 2

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "'011827"
$ws.Range("B16").ClearFormats()
$ws.Range("C16").Value = '汇添富健康生活一年持有期混合C'
$ws.Range("D16").Value = "'3.55"
$ws.Range("D16").ClearFormats()
$ws.Range("E16").Value = "'88.44"
$ws.Range("E16").ClearFormats()
$ws.Range("F16").Value = "'3.38"
$ws.Range("F16").ClearFormats()
$ws.Range("G16").Value = "'0.1200"
$ws.Range("G16").ClearFormats()
$ws.Range("H16").Value = 9

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "'011650"
$ws.Range("B17").ClearFormats()
$ws.Range("C17").Value = '易方达逆向投资混合C'
$ws.Range("D17").Value = "'1.96"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'85.02"
$ws.Range("E17").ClearFormats()
$ws.Range("F17").Value = "'5.82"
$ws.Range("F17").ClearFormats()
$ws.Range("G17").Value = "'0.1141"
$ws.Range("G17").ClearFormats()
$ws.Range("H17").Value = 1

$ws.Range("A18").Value = 16
$ws.Range("B18").Value = "'011845"
$ws.Range("B18").ClearFormats()
$ws.Range("C18").Value = '博时周期优选混合型证券投资基金A'
$ws.Range("D18").Value = "'2.36"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'86.75"
$ws.Range("E18").ClearFormats()
$ws.Range("F18").Value = "'2.90"
$ws.Range("F18").ClearFormats()
$ws.Range("G18").Value = "'0.0684"
$ws.Range("G18").ClearFormats()
$ws.Range("H18").Value = 8

$ws.Range("A19").Value = 17
$ws.Range("B19").Value = "'009189"
$ws.Range("B19").ClearFormats()
$ws.Range("C19").Value = '华宝成长策略混合'
$ws.Range("D19").Value = "'1.47"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'83.71"
$ws.Range("E19").ClearFormats()
$ws.Range("F19").Value = "'4.57"
$ws.Range("F19").ClearFormats()
$ws.Range("G19").Value = "'0.0672"
$ws.Range("G19").ClearFormats()
$ws.Range("H19").Value = 3

$ws.Range("A20").Value = 18
$ws.Range("B20").Value = "'000867"
$ws.Range("B20").ClearFormats()
$ws.Range("C20").Value = '华宝品质生活股票'
$ws.Range("D20").Value = "'0.81"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'87.64"
$ws.Range("E20").ClearFormats()
$ws.Range("F20").Value = "'5.23"
$ws.Range("F20").ClearFormats()
$ws.Range("G20").Value = "'0.0424"
$ws.Range("G20").ClearFormats()
$ws.Range("H20").Value = 3

$ws.Range("A21").Value = 19
$ws.Range("B21").Value = "'010390"
$ws.Range("B21").ClearFormats()
$ws.Range("C21").Value = '易方达科益混合C'
$ws.Range("D21").Value = "'0.29"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'92.94"
$ws.Range("E21").ClearFormats()
$ws.Range("F21").Value = "'8.85"
$ws.Range("F21").ClearFormats()
$ws.Range("G21").Value = "'0.0257"
$ws.Range("G21").ClearFormats()
$ws.Range("H21").Value = 2

$ws.Range("A22").Value = 20
$ws.Range("B22").Value = "'002319"
$ws.Range("B22").ClearFormats()
$ws.Range("C22").Value = '大成一带一路灵活配置混合'
$ws.Range("D22").Value = "'0.50"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'89.30"
$ws.Range("E22").ClearFormats()
$ws.Range("F22").Value = "'5.13"
$ws.Range("F22").ClearFormats()
$ws.Range("G22").Value = "'0.0256"
$ws.Range("G22").ClearFormats()
$ws.Range("H22").Value = 7

$ws.Range("A23").Value = 21
$ws.Range("B23").Value = "'011846"
$ws.Range("B23").ClearFormats()
$ws.Range("C23").Value = '博时周期优选混合型证券投资基金C'
$ws.Range("D23").Value = "'0.07"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'86.75"
$ws.Range("E23").ClearFormats()
$ws.Range("F23").Value = "'2.90"
$ws.Range("F23").ClearFormats()
$ws.Range("G23").Value = "'0.0020"
$ws.Range("G23").ClearFormats()
$ws.Range("H23").Value = 8

# ---------------------------------------------------------------------------
# 3) Update "总计": push the three existing quarters down one row and add the
#    new 2022-Q1 summary at the top (row 2).
# ---------------------------------------------------------------------------
$total = $wb.Worksheets.Item(5)

# extend the styled column-A border down to the new row 5
$total.Range("A4").Copy()
$total.Range("A5").PasteSpecial(-4122)

$total.Range("B5").Value = $total.Range("B4").Value2
$total.Range("C5").Value = $total.Range("C4").Value2
$total.Range("D5").Value = $total.Range("D4").Value2
$total.Range("A5").Value = 3

$total.Range("B4").Value = $total.Range("B3").Value2
$total.Range("C4").Value = $total.Range("C3").Value2
$total.Range("D4").Value = $total.Range("D3").Value2

$total.Range("B3").Value = $total.Range("B2").Value2
$total.Range("C3").Value = $total.Range("C2").Value2
$total.Range("D3").Value = $total.Range("D2").Value2

$total.Range("B2").Value = "2022-Q1"
$total.Range("C2").Value = 22
$total.Range("D2").Value = 6.93

# ---------------------------------------------------------------------------
# 4) Restore the originally active tab ("2021-Q2"), since Worksheet.Copy()
#    makes the new sheet active.
# ---------------------------------------------------------------------------
$wb.Worksheets.Item(1).Activate()
